# Automatische test-sync: 2025-06-29 15:15:50
# Appends a new test-mail log entry (row 31) to the "Logs" sheet and
# updates the corresponding category count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

$newRow = 31

$logs.Cells.Item($newRow, 1).Value = "Kun je dit product voor mij bestellen?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #16: Kun je dit product voor mij bestellen?"
$logs.Cells.Item($newRow, 4).Value = "Bestelling / Levering"
$logs.Cells.Item($newRow, 5).Value = "Geachte klant,`nDank u wel voor uw interesse in ons product. Helaas kunnen wij op basis van deze e-mail geen bestelling voor u plaatsen. U kunt echter onze website bezoeken en het product zelf bestellen. Mocht u verdere vragen hebben of hulp nodig hebben bij het plaatsen van een bestelling, dan helpen wij u graag verder.`nMet vriendelijke groet,`n[Naam Bedrijf] E-mailassistent"
$logs.Cells.Item($newRow, 6).Value = "2025-06-29 15:15:43"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"

# The new cell in column E contains embedded newlines; clear the
# auto-computed "custom height" that the runtime stamps on such rows so
# row 31 stays un-pinned, same as every other data row in this sheet.
$logs.Rows.Item($newRow).AutoFit()

# Extend the conditional-formatting rule ranges (Category / Beantwoord /
# Handmatig opvolgen / Automatisch afgehandeld) so the new row inherits
# the same colouring as the rest of the log.
$oldRanges = @("D2:D30", "G2:G30", "H2:H30", "I2:I30")
$newRanges = @("D2:D31", "G2:G31", "H2:H31", "I2:I31")

for ($i = 0; $i -lt $oldRanges.Length; $i++) {
    $fcs = $logs.Range($oldRanges[$i]).FormatConditions
    for ($j = 1; $j -le $fcs.Count; $j++) {
        $fcs.Item($j).ModifyAppliesToRange($logs.Range($newRanges[$i]))
    }
}

# Bump the "Bestelling / Levering" tally on the Dashboard sheet to
# account for the new log entry.
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 9
